$d = $word.ActiveDocument

# 1. Split "It should add " into "It should a" + bookmark(_GoBack) + "dd "
#    Find the paragraph text first so the bookmark can be re-inserted mid-run.
$r = $d.Content
$r.Find.Execute("It should a", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$r.Collapse(0)
$d.Bookmarks.Add("_GoBack", $r)

# 2. Remove the two trailing empty paragraphs that used to hold the bookmark.
$end = $d.Content.End
$r2 = $d.Range($end - 2, $end)
$r2.Delete()
